# Switch to English federal state names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "Bayern"              = "Bavaria"
    "Niedersachsen"       = "Lower Saxony"
    "Nordrhein-Westfalen" = "North Rhine-Westphalia"
    "Rheinland-Pfalz"     = "Rhineland-Palatinate"
    "Sachsen"             = "Saxony"
    "Sachsen-Anhalt"      = "Saxony-Anhalt"
    "Thüringen"           = "Thuringia"
}

$used = $ws.UsedRange
foreach ($row in 1..$used.Rows.Count) {
    $cell = $ws.Cells.Item($row, 1)
    $val = $cell.Text
    if ($null -ne $val -and $replacements.ContainsKey([string]$val)) {
        $cell.Value = $replacements[[string]$val]
    }
}

# Match the final cursor/selection position recorded in the saved file.
$ws.Range("H20").Select() | Out-Null
